# Contingencies with rene fine
# Insert two new line rows ("line7", "line8") into the lines/extractions
# table (pushing the existing extr1..extr8 rows down by two), and
# refresh the from_bus / to_bus / in_service values for every
# contingency row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right after "line6" (row 7), before the old
# "extr1" row (old row 8) - this shifts the old extr1..extr8 rows
# (old rows 8-15) down to rows 10-17, matching the row where Excel
# would place them after sorting the new lines/extr split.
$ws.Rows("8:9").Insert(-4121) | Out-Null

# Give the two new "A" index cells the same look (bold, bordered,
# centered) as the rest of the index column before filling them in.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A8:A9").PasteSpecial(-4122) | Out-Null

# ---- New row 8: line7 ----
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# ---- New row 9: line8 ----
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# ---- Refreshed contingency rows (old extr1..extr8, now rows 10-17) ----
# The "A" index column needs to be bumped by 2 as well, since these
# rows moved down two positions in the table.
$ws.Range("A10").Value = 8
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

$ws.Range("A11").Value = 9
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

$ws.Range("A12").Value = 10
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = $true

$ws.Range("A13").Value = 11
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $false

$ws.Range("A14").Value = 12
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $false

$ws.Range("A15").Value = 13
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $false

$ws.Range("A16").Value = 14
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

$ws.Range("A17").Value = 15
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true
